$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 212, shifting the existing rows
# 212-231 down to 214-233 (weekly data refresh: two fresh rows of data
# land at the top of this block).
$ws.Range("A212:A213").EntireRow.Insert()

# Row 212 - new weekly entry
$ws.Range("A212").Value = 9
$ws.Range("B212").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C212").Value = "Metropolitana"
$ws.Range("D212").Value = 44461
$ws.Range("E212").Value = 13
$ws.Range("F212").Value = 100112012
$ws.Range("G212").Value = "Espinaca"
$ws.Range("H212").Value = "Sin especificar"
$ws.Range("I212").Value = "Primera"
$ws.Range("J212").Value = 133
$ws.Range("K212").Value = 6000
$ws.Range("L212").Value = 7000
$ws.Range("M212").Value = 6504
$ws.Range("N212").Value = '$/cuna 10 kilos'
$ws.Range("O212").Value = "Provincia de Chacabuco"
$ws.Range("P212").Value = 650
$ws.Range("Q212").Value = 10
$ws.Range("R212").Value = "Hortaliza"

# Row 213 - new weekly entry
$ws.Range("A213").Value = 9
$ws.Range("B213").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C213").Value = "Metropolitana"
$ws.Range("D213").Value = 44461
$ws.Range("E213").Value = 13
$ws.Range("F213").Value = 100112012
$ws.Range("G213").Value = "Espinaca"
$ws.Range("H213").Value = "Sin especificar"
$ws.Range("I213").Value = "Segunda"
$ws.Range("J213").Value = 61
$ws.Range("K213").Value = 5000
$ws.Range("L213").Value = 5000
$ws.Range("M213").Value = 5000
$ws.Range("N213").Value = '$/cuna 10 kilos'
$ws.Range("O213").Value = "Provincia de Chacabuco"
$ws.Range("P213").Value = 500
$ws.Range("Q213").Value = 10
$ws.Range("R213").Value = "Hortaliza"

# Apply the same date-display style used by the rest of column D (style
# index 2 in styles.xml -> YYYY-MM-DD HH:MM:SS numFmt) to the two new
# date cells, matching their neighbours.
$ws.Range("D214").Copy()
$ws.Range("D212:D213").PasteSpecial(-4122)
$ws.Range("D212").Value = 44461
$ws.Range("D213").Value = 44461
